$d = $word.ActiveDocument

$d.Content.Find.Execute("** 11:45 AM (Estimado)", $true, $false, $false, $false, $false, $true, 1, $false, "** No especificado en transcripción", 2)

$d.Content.Find.Execute("COORDINACIÓN ACADÉMICA: JOHON FREDY SANABRIA MUÑOZ", $true, $false, $false, $false, $false, $true, 1, $false, "COORDINACIÓN ACADÉMICA: Ingeniero John Freddy Sanabria Muñoz", 2)

$d.Content.Find.Execute("BIENESTAR DEL APRENDIZ: Dra. Elizabeth (No apellido especificado en transcripción)", $true, $false, $false, $false, $false, $true, 1, $false, "BIENESTAR DEL APRENDIZ: Doctora Elizabeth", 2)

$d.Content.Find.Execute("INSTRUCTORES: Viviana Barrera (Profe de ICB)", $true, $false, $false, $false, $false, $true, 1, $false, "INSTRUCTORES: Viviana Barrera", 2)

$d.Content.Find.Execute("REPRESENTANTE DE CENTRO: Karen Andrea García (Líder de etapas productivas), Diana Alicia Alfonso (Líder de contratos)", $true, $false, $false, $false, $false, $true, 1, $false, "REPRESENTANTE DE CENTRO: No especificado en transcripción", 2)

$d.Content.Find.Execute("VOCERO: Alejandra Sely (Jefe de Talento Humano - Hotel Sonesta)", $true, $false, $false, $false, $false, $true, 1, $false, "VOCERO: No especificado en transcripción", 2)
